$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text looks numeric: force Text format so Excel
# does not silently convert the string to a number (losing formatting
# like trailing zeros), matching the original inline-string cells.
$textCells = @("D5","D6","D9","D10","D11","D12","D13","D14","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D30","D31","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D46","D47","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = '552.07'
$ws.Range("D6").Value = '136.12'
$ws.Range("D9").Value = '0.491'
$ws.Range("D10").Value = '6.58'
$ws.Range("D11").Value = '0.157'
$ws.Range("D12").Value = '0.451'
$ws.Range("D13").Value = '34.74'
$ws.Range("D14").Value = '0.0000215'
$ws.Range("D19").Value = '500.11'
$ws.Range("D20").Value = '6.61'
$ws.Range("D21").Value = '13.40'
$ws.Range("D22").Value = '0.699'
$ws.Range("D23").Value = '7.20'
$ws.Range("D24").Value = '77.06'
$ws.Range("D25").Value = '12.16'
$ws.Range("D26").Value = '0.999'
$ws.Range("D27").Value = '2.74'
$ws.Range("D28").Value = '8.08'
$ws.Range("D30").Value = '1.96'
$ws.Range("D31").Value = '26.08'
$ws.Range("D33").Value = '2.48'
$ws.Range("D34").Value = '58.51'
$ws.Range("D35").Value = '524.87'
$ws.Range("D36").Value = '5.83'
$ws.Range("D37").Value = '5.13'
$ws.Range("D38").Value = '0.0408'
$ws.Range("D40").Value = '0.0784'
$ws.Range("D41").Value = '0.120'
$ws.Range("D42").Value = '8.00'
$ws.Range("D43").Value = '2.58'
$ws.Range("D44").Value = '0.251'
$ws.Range("D46").Value = '2.03'
$ws.Range("D47").Value = '121.25'
$ws.Range("D48").Value = '0.105'
$ws.Range("D49").Value = '23.49'

# Restore the default (unstyled) cell style now that the text value is set,
# so the style index matches the rest of the untouched price cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining plain text / percentage / label updates
$ws.Range("D2").Value = '62.973.52'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '3.078.98'
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("E6").Value = '  -3.09%  '
$ws.Range("D8").Value = '3.060.13'
$ws.Range("E8").Value = '  +0.43%  '
$ws.Range("E9").Value = '  +1.12%  '
$ws.Range("E10").Value = '  +1.79%  '
$ws.Range("E11").Value = '  +4.92%  '
$ws.Range("E12").Value = '  +1.63%  '
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("D15").Value = '3.568.33'
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = '63.116.70'
$ws.Range("E16").Value = '  -0.73%  '
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = '3.089.84'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("E19").Value = '  +3.68%  '
$ws.Range("E20").Value = '  +1.77%  '
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("E22").Value = '  +4.00%  '
$ws.Range("E23").Value = '  +1.60%  '
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("E27").Value = '  +2.46%  '
$ws.Range("E28").Value = '  -0.45%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  -3.96%  '
$ws.Range("E31").Value = '  +2.56%  '
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("E33").Value = '  -3.65%  '
$ws.Range("E34").Value = '  +12.95%  '
$ws.Range("E35").Value = '  -8.63%  '
$ws.Range("E36").Value = '  +0.83%  '
$ws.Range("E37").Value = '  -2.30%  '
$ws.Range("E38").Value = '  +3.51%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '3.036.83'
$ws.Range("E39").Value = '  +1.54%  '
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("E41").Value = '  +3.85%  '
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("E43").Value = '  -7.59%  '
$ws.Range("E44").Value = '  +4.57%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("E47").Value = '  +3.80%  '
$ws.Range("E48").Value = '  -0.30%  '
$ws.Range("E49").Value = '  -3.99%  '
$ws.Range("D50").Value = '0.0₃0494'
$ws.Range("E50").Value = '  -3.75%  '
$ws.Range("E51").Value = '  +67.81%  '
